# Excel COM-interop edit script
# Commit: "Update countries & provincias Spain"
# - Refreshes the COVID-19 "paises" dataset with a newer snapshot (01:26 update).
# - Chequia (Czechia) overtakes Canada in the total-cases ranking, so rows 31/32
#   swap which country label they carry, each keeping its own up-to-date numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 01:26"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8578201
$ws.Range("C4").Value = 57045
$ws.Range("D4").Value = 5585140
$ws.Range("E4").Value = 2765760
$ws.Range("G4").Value = 1117
$ws.Range("H4").Value = 227301

# Row 9: Argentina
$ws.Range("B9").Value = 1037325
$ws.Range("C9").Value = 18326
$ws.Range("D9").Value = 840520
$ws.Range("E9").Value = 169286
$ws.Range("G9").Value = 419
$ws.Range("H9").Value = 27519

# Row 12: Peru
$ws.Range("B12").Value = 876885
$ws.Range("C12").Value = 2767
$ws.Range("D12").Value = 792892
$ws.Range("E12").Value = 50056
$ws.Range("G12").Value = 62
$ws.Range("H12").Value = 33937

# Row 21: Alemania
$ws.Range("B21").Value = 391355
$ws.Range("C21").Value = 10457
$ws.Range("E21").Value = 79256

# Row 30: Belgica
$ws.Range("D30").Value = 21477
$ws.Range("E30").Value = 208193

# Row 31: Chequia
$ws.Range("A31").Value = "Chequia"
$ws.Range("B31").Value = 208915
$ws.Range("C31").Value = 14969
$ws.Range("D31").Value = 83136
$ws.Range("E31").Value = 124040
$ws.Range("G31").Value = 120
$ws.Range("H31").Value = 1739

# Row 32: Canada
$ws.Range("A32").Value = "Canada"
$ws.Range("B32").Value = 205749
$ws.Range("C32").Value = 2061
$ws.Range("D32").Value = 173392
$ws.Range("E32").Value = 22533
$ws.Range("G32").Value = 30
$ws.Range("H32").Value = 9824

# Row 36: Ecuador
$ws.Range("B36").Value = 155625
$ws.Range("C36").Value = 1510
$ws.Range("E36").Value = 8985
$ws.Range("G36").Value = 49
$ws.Range("H36").Value = 12453

# Row 40: Panama
$ws.Range("B40").Value = 126435
$ws.Range("C40").Value = 696
$ws.Range("D40").Value = 102725
$ws.Range("E40").Value = 21113
$ws.Range("G40").Value = 12
$ws.Range("H40").Value = 2597

# Row 58: Barein
$ws.Range("B58").Value = 78907
$ws.Range("C58").Value = 374
$ws.Range("D58").Value = 75424
$ws.Range("E58").Value = 3175
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 308

# Row 64: Nigeria
$ws.Range("B64").Value = 61667
$ws.Range("C64").Value = 37
$ws.Range("D64").Value = 56880
$ws.Range("E64").Value = 3662

# Row 66: Paraguay
$ws.Range("B66").Value = 56819
$ws.Range("C66").Value = 746
$ws.Range("D66").Value = 37673
$ws.Range("E66").Value = 17896
$ws.Range("G66").Value = 19
$ws.Range("H66").Value = 1250

# Row 96: Noruega
$ws.Range("B96").Value = 16964
$ws.Range("C96").Value = 193
$ws.Range("E96").Value = 4822

# Row 114: Gabon
$ws.Range("B114").Value = 8901
$ws.Range("C114").Value = 17
$ws.Range("D114").Value = 8479
$ws.Range("E114").Value = 368

# Row 155: Uruguay
$ws.Range("B155").Value = 2663
$ws.Range("C155").Value = 40
$ws.Range("D155").Value = 2172
$ws.Range("E155").Value = 438
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 53
